# Update the "K" column (column G) with recalculated strikeout values.
# Commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 1;  3 = 0;  4 = 3;  5 = 3;  6 = 2;  7 = 1;  8 = 2;  9 = 1;  10 = 2;
    11 = 2;  12 = 7;  13 = 4;  14 = 2;  15 = 0;  16 = 1;  17 = 3;  18 = 1;  19 = 1;  20 = 0;
    21 = 2;  22 = 4;  23 = 3;  24 = 3;  25 = 0;  26 = 0;  27 = 3;  28 = 4;  29 = 2;  30 = 0;
    31 = 1;  32 = 2;  33 = 2;  34 = 0;  35 = 1;  36 = 3;  37 = 3;  38 = 2;  39 = 1;  40 = 1;
    41 = 7;  42 = 2;  43 = 1;  44 = 2;  45 = 1;  46 = 3;  47 = 2;  48 = 4;  49 = 0;  50 = 2;
    51 = 0;  52 = 3;  53 = 0;  54 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
